$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (buses 0-23), columns B,C,D,E,F,I,J,K,L,M,N
# (column G stays 1, column H / A are untouched per the diff)
$data = @{
    2 = @(1.02, 1.017946612101155, 1.020591546546963, 1.019281350786939, 1.016317380383342, 1.027641855197841, 1.023157837374231, 1.02343137844637, 1.022125063242104, 1.019169911772201, 1.011820112493132)
    3 = @(1.02, 1.018911443448249, 1.021468284749887, 1.020099142501327, 1.017924006625639, 1.027773324297557, 1.023758696230691, 1.0241141953398, 1.022748813804788, 1.020579675058588, 1.012018805161629)
    4 = @(1.02, 1.019535169988972, 1.02203530340856, 1.02062821644213, 1.018962888425043, 1.027855960205031, 1.024146341619783, 1.024555104052296, 1.023151686628898, 1.021490715763239, 1.012146968577176)
    5 = @(1.02, 1.019797245339543, 1.022273608731252, 1.020850616927766, 1.019399470051224, 1.027890117375177, 1.024309032487328, 1.024740242067041, 1.023320878307281, 1.021873442507086, 1.012200751776611)
    6 = @(1.02, 1.019841240790148, 1.022313617178314, 1.020887957640658, 1.019472764578137, 1.027895818316864, 1.02433633286348, 1.024771314656051, 1.023349276010882, 1.021937688056155, 1.012209776542498)
    7 = @(1.02, 1.01953867239936, 1.022038487928365, 1.020631188255793, 1.018968722687749, 1.027856418906009, 1.024148516584504, 1.024557578739886, 1.023153948068177, 1.021495830849676, 1.012147687610869)
    8 = @(1.02, 1.018272802563717, 1.020887904297189, 1.01955774626946, 1.016860498677262, 1.027686789405563, 1.023361138758335, 1.023662330071729, 1.022336014973936, 1.019646593493263, 1.011887345283746)
    9 = @(1.02, 1.016037683233192, 1.018858209290455, 1.017665506404562, 1.013139793280231, 1.027369260671296, 1.021964858261938, 1.022077738853534, 1.020889077026778, 1.016378792874103, 1.011425492933618)
    10 = @(1.02, 1.014544544778855, 1.017503580100079, 1.016403548920754, 1.010655057955061, 1.027145076659251, 1.021028058357137, 1.021016592543038, 1.019920656162811, 1.014193724203689, 1.011115506479536)
    11 = @(1.02, 1.01389726435544, 1.016916652275561, 1.01585699714181, 1.009578038339762, 1.027045041800716, 1.020620998626921, 1.0205559732612, 1.019500416424995, 1.013245938085777, 1.010980783290468)
    12 = @(1.02, 1.013656723213948, 1.016698586015191, 1.015653965931143, 1.009177810256777, 1.027007439644481, 1.020469584781074, 1.02038470761593, 1.019344184058745, 1.012893635433705, 1.010930666291643)
    13 = @(1.02, 1.013708325205903, 1.016745364469498, 1.015697517611389, 1.009263668586953, 1.027015525557613, 1.020502073245632, 1.020421452430259, 1.01937770261936, 1.012969217092065, 1.01094141994551)
    14 = @(1.02, 1.01387738344758, 1.016898627977426, 1.015840214877374, 1.009544958978031, 1.027041942666494, 1.020608487074453, 1.020541819883452, 1.019487504997787, 1.013216821832638, 1.010976642130589)
    15 = @(1.02, 1.013981530890696, 1.016993051365578, 1.015928133043249, 1.009718247840874, 1.027058160190618, 1.020674023808273, 1.020615959552417, 1.019555139753806, 1.013369345598905, 1.010998333764294)
    16 = @(1.02, 1.014587486942901, 1.01754252485642, 1.016439819309902, 1.010726512084748, 1.027151653261673, 1.021055043626122, 1.021047138355012, 1.019948526950331, 1.014256590582663, 1.011124437130864)
    17 = @(1.02, 1.014967387947245, 1.017887097226842, 1.016760755476086, 1.011358666146598, 1.027209506214127, 1.021293666972969, 1.021317300977791, 1.020195045256383, 1.014812692908421, 1.011203405282479)
    18 = @(1.02, 1.01518890628513, 1.018088045323701, 1.01694794094346, 1.011727284074118, 1.027242965088218, 1.021432715038766, 1.021474772697341, 1.020338747713599, 1.015136900498657, 1.011249418124738)
    19 = @(1.02, 1.015264426216426, 1.018156557483129, 1.017011764460053, 1.011852955377428, 1.02725432522084, 1.021480103631622, 1.021528447913865, 1.020387731706261, 1.015247420400686, 1.011265099186742)
    20 = @(1.02, 1.014926635566036, 1.017850131509405, 1.016726323218781, 1.011290853105343, 1.027203328697213, 1.021268079109427, 1.021288326421456, 1.020168605227917, 1.014753044729429, 1.011194937707215)
    21 = @(1.02, 1.013827603118836, 1.016853497219436, 1.015798194574701, 1.009462130864243, 1.027034175762209, 1.020577156758641, 1.020506379396638, 1.019455174708136, 1.013143915437807, 1.010966272141175)
    22 = @(1.02, 1.013135947163069, 1.016226554497364, 1.015214541668285, 1.00831132393323, 1.026925249470177, 1.020141509562248, 1.020013747883694, 1.019005822026281, 1.012130726634682, 1.010822068136094)
    23 = @(1.02, 1.013502669067791, 1.016558939084097, 1.015523956797832, 1.008921487468315, 1.026983237196416, 1.020372571805722, 1.020274995210258, 1.019244107440395, 1.012667978355892, 1.010898554488137)
    24 = @(1.02, 1.014945050033463, 1.017866834834261, 1.016741881708056, 1.011321495227971, 1.02720612093412, 1.021279641586271, 1.021301419114558, 1.020180552608831, 1.01477999765897, 1.011198763988056)
    25 = @(1.02, 1.016616051540808, 1.019383197764602, 1.018154777887316, 1.014102406914045, 1.027453553812817, 1.022326877056411, 1.022488230238871, 1.02126381360748, 1.017224724135029, 1.011545260041226)
}

$cols = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
